$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "252.01") are not coerced into floating point numbers, then restore
# the default "Normal" style so no stray style index is left referenced.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.002.12"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "2.291.16"
$ws.Range("E3").Value = "  +1.81%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "252.01"
$ws.Range("E5").Value = "  -0.51%  "

$ws.Range("D6").Value = "0.642"
$ws.Range("E6").Value = "  +1.11%  "

$ws.Range("D7").Value = "73.94"
$ws.Range("E7").Value = "  +4.60%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "39.24"
$ws.Range("E10").Value = "  -4.54%  "

$ws.Range("D11").Value = "0.0981"
$ws.Range("E11").Value = "  +3.32%  "

$ws.Range("D12").Value = "59.10"
$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("D13").Value = "7.44"
$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("E14").Value = "  +1.62%  "

$ws.Range("D15").Value = "2.634.00"
$ws.Range("E15").Value = "  +1.71%  "

$ws.Range("D16").Value = "15.38"
$ws.Range("E16").Value = "  +3.66%  "

$ws.Range("D17").Value = "0.876"
$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("D18").Value = "2.289.46"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").Value = "42.884.29"
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("E20").Value = "  +2.96%  "

$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("D22").Value = "72.70"
$ws.Range("E22").Value = "  -0.44%  "

$ws.Range("D23").Value = "238.15"
$ws.Range("E23").Value = "  +1.07%  "

$ws.Range("D24").Value = "2.24"
$ws.Range("E24").Value = "  +5.27%  "

$ws.Range("E25").Value = "  -2.49%  "

$ws.Range("D26").Value = "11.63"
$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("D28").Value = "2.42"
$ws.Range("E28").Value = "  -0.78%  "

$ws.Range("E29").Value = "  -1.08%  "

$ws.Range("E30").Value = "  -3.20%  "

$ws.Range("D31").Value = "167.13"
$ws.Range("E31").Value = "  -0.43%  "

$ws.Range("D32").Value = "21.08"

$ws.Range("E33").Value = "  +5.32%  "

$ws.Range("D34").Value = "0.127"
$ws.Range("E34").Value = "  +2.52%  "

$ws.Range("D35").Value = "0.0829"
$ws.Range("E35").Value = "  +5.12%  "

$ws.Range("D36").Value = "31.04"
$ws.Range("E36").Value = "  +10.94%  "

$ws.Range("E37").Value = "  +2.62%  "

$ws.Range("E38").Value = "  +10.01%  "

$ws.Range("D39").Value = "4.78"
$ws.Range("E39").Value = "  +1.36%  "

$ws.Range("E40").Value = "  -3.30%  "

$ws.Range("D41").Value = "14.33"
$ws.Range("E41").Value = "  +13.50%  "

$ws.Range("E42").Value = "  +2.73%  "

$ws.Range("D43").Value = "5.94"
$ws.Range("E43").Value = "  +1.78%  "

$ws.Range("E44").Value = "  +6.01%  "

$ws.Range("E45").Value = "  +4.38%  "

$ws.Range("D46").Value = "61.86"
$ws.Range("E46").Value = "  -3.62%  "

$ws.Range("D47").Value = "4.90"
$ws.Range("E47").Value = "  -1.66%  "

$ws.Range("D48").Value = "0.104"
$ws.Range("E48").Value = "  +1.20%  "

$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("E50").Value = "  -1.82%  "

$ws.Range("D51").Value = "100.26"
$ws.Range("E51").Value = "  +6.06%  "

# Restore default styling on column D (removes the temporary text format).
$dRange.Style = "Normal"
